$wb = $excel.ActiveWorkbook

# The edited sheet is "dotnet" (first sheet, internally sheet1.xml)
$ws = $wb.Worksheets.Item("dotnet")
$ws.Activate()

# Update Published_File_Path (C2) and Destination_Path (D2) values
$ws.Range("C2").Value = "D:\\Git\\demo\\petmatrix-backend-kku\\PetMatrix.API\\bin\\Release\\netcoreapp2.1\\publish"
$ws.Range("D2").Value = "D:\Publish Files\PetMatrix\PetMatrixBackend"

# Widen column D to fit the new content (closest attainable width to 47.140625)
$ws.Columns.Item(4).ColumnWidth = 46.3

# Update the selected cell to D9
$ws.Range("D9").Select()
